$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "SCD0014"

# Update the TC_ID cell
$ws.Range("B2").Value = "SCD0014-001"

$rng = $ws.Range("A1:U6")
$rng.HorizontalAlignment = -4131
